$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell B2 value from 2 to 2030
$ws.Range("B2").Value = 2030

# Update the active selection from C4 to B3
$ws.Range("B3").Select()
